$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.859.84'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.881.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.41%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.53'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4674'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3939'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07935'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.42'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.920.70'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.752'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.013'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06993'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.65'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001013'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.01'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.859.22'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.367'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.13'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.124'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.122.66'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.56'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.45'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.765'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.009'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.05'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09400'
$ws.Range('D31').Style = "Normal"
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('E34').Value = '  +3.25%  '
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05926'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02123'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.152'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  +4.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5723'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.01'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1796'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07265'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.80'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5347'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.152'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.135'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.854'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '114.41'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.374'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.33%  '
